# Add cm015 to syllabus:
#  - mark the "link_it" column (C16) as TRUE now that the class materials are linked
#  - shorten the topic (D16) from "Hypothesis testing and p-values" to "Hypothesis testing"
#  - leave the cursor on the next row (D17), matching where the user clicked next

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D16").Value = "Hypothesis testing"
$ws.Range("C16").Value = $true

$ws.Range("D17").Select()
